$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.043260161776686
$ws.Cells.Item(2, 4).Value = 1.046428770128438
$ws.Cells.Item(2, 5).Value = 1.041230987379512
$ws.Cells.Item(2, 6).Value = 1.054716789940024
$ws.Cells.Item(2, 9).Value = 1.035695167092103
$ws.Cells.Item(2, 10).Value = 1.048331565347354
$ws.Cells.Item(2, 11).Value = 1.049194232275643
$ws.Cells.Item(2, 12).Value = 1.044011088834224
$ws.Cells.Item(2, 13).Value = 1.057459232750838

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.0447761258573
$ws.Cells.Item(3, 4).Value = 1.04757475203088
$ws.Cells.Item(3, 5).Value = 1.042538705071643
$ws.Cells.Item(3, 6).Value = 1.056094196138114
$ws.Cells.Item(3, 9).Value = 1.035997814210861
$ws.Cells.Item(3, 10).Value = 1.049491494871725
$ws.Cells.Item(3, 11).Value = 1.0501509862116
$ws.Cells.Item(3, 12).Value = 1.045128086639668
$ws.Cells.Item(3, 13).Value = 1.058648496899756

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.045755743955565
$ws.Cells.Item(4, 4).Value = 1.04831492629831
$ws.Cells.Item(4, 5).Value = 1.043383953425448
$ws.Cells.Item(4, 6).Value = 1.05698447058277
$ws.Cells.Item(4, 9).Value = 1.036191685584827
$ws.Cells.Item(4, 10).Value = 1.050240359924659
$ws.Cells.Item(4, 11).Value = 1.050768141607326
$ws.Cells.Item(4, 12).Value = 1.045849387169562
$ws.Cells.Item(4, 13).Value = 1.059416500872544

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.046167268682975
$ws.Cells.Item(5, 4).Value = 1.048625776260938
$ws.Cells.Item(5, 5).Value = 1.043739078018836
$ws.Cells.Item(5, 6).Value = 1.057358508280055
$ws.Cells.Item(5, 9).Value = 1.03627272108411
$ws.Cells.Item(5, 10).Value = 1.050554784244984
$ws.Cells.Item(5, 11).Value = 1.051027136324884
$ws.Cells.Item(5, 12).Value = 1.046152274529214
$ws.Cells.Item(5, 13).Value = 1.059739008672601

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.046236347723033
$ws.Cells.Item(6, 4).Value = 1.048677950744641
$ws.Cells.Item(6, 5).Value = 1.043798692384657
$ws.Cells.Item(6, 6).Value = 1.057421297314222
$ws.Cells.Item(6, 9).Value = 1.036286299901193
$ws.Cells.Item(6, 10).Value = 1.050607554254348
$ws.Cells.Item(6, 11).Value = 1.051070595977686
$ws.Cells.Item(6, 12).Value = 1.04620311042239
$ws.Cells.Item(6, 13).Value = 1.059793138125189

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.04576124396519
$ws.Cells.Item(7, 4).Value = 1.04831908113977
$ws.Cells.Item(7, 5).Value = 1.04338869946682
$ws.Cells.Item(7, 6).Value = 1.056989469404074
$ws.Cells.Item(7, 9).Value = 1.036192770222228
$ws.Cells.Item(7, 10).Value = 1.050244562837297
$ws.Cells.Item(7, 11).Value = 1.050771604099502
$ws.Cells.Item(7, 12).Value = 1.045853435724152
$ws.Cells.Item(7, 13).Value = 1.059420811650039

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.043772765054363
$ws.Cells.Item(8, 4).Value = 1.046816342888006
$ws.Cells.Item(8, 5).Value = 1.041673133030758
$ws.Cells.Item(8, 6).Value = 1.05518250113821
$ws.Cells.Item(8, 9).Value = 1.035797855204747
$ws.Cells.Item(8, 10).Value = 1.048723921447507
$ws.Cells.Item(8, 11).Value = 1.049517973595136
$ws.Cells.Item(8, 12).Value = 1.044388891197042
$ws.Cells.Item(8, 13).Value = 1.057861469891345

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.040258430923221
$ws.Cells.Item(9, 4).Value = 1.044157762431224
$ws.Cells.Item(9, 5).Value = 1.038642700955887
$ws.Cells.Item(9, 6).Value = 1.051990485707186
$ws.Cells.Item(9, 9).Value = 1.035086870374436
$ws.Cells.Item(9, 10).Value = 1.046031198798407
$ws.Cells.Item(9, 11).Value = 1.047293956184705
$ws.Cells.Item(9, 12).Value = 1.041796681934693
$ws.Cells.Item(9, 13).Value = 1.055101767873639

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.037908096796548
$ws.Cells.Item(10, 4).Value = 1.042377985299028
$ws.Cells.Item(10, 5).Value = 1.036617105066968
$ws.Cells.Item(10, 6).Value = 1.049856805095535
$ws.Cells.Item(10, 9).Value = 1.034602631456029
$ws.Cells.Item(10, 10).Value = 1.044226867822961
$ws.Cells.Item(10, 11).Value = 1.045800949793095
$ws.Cells.Item(10, 12).Value = 1.040060500245659
$ws.Cells.Item(10, 13).Value = 1.053253614676644

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.036888506212694
$ws.Cells.Item(11, 4).Value = 1.04160550449564
$ws.Cells.Item(11, 5).Value = 1.035738663799029
$ws.Cells.Item(11, 6).Value = 1.048931473296447
$ws.Cells.Item(11, 9).Value = 1.034390496916137
$ws.Cells.Item(11, 10).Value = 1.04344331925545
$ws.Cells.Item(11, 11).Value = 1.045151951151662
$ws.Cells.Item(11, 12).Value = 1.039306736441616
$ws.Cells.Item(11, 13).Value = 1.052451290549654

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.036509493311694
$ws.Cells.Item(12, 4).Value = 1.041318290781311
$ws.Cells.Item(12, 5).Value = 1.035412162834608
$ws.Cells.Item(12, 6).Value = 1.048587541471988
$ws.Cells.Item(12, 9).Value = 1.034311329702884
$ws.Cells.Item(12, 10).Value = 1.043151928106744
$ws.Cells.Item(12, 11).Value = 1.044910500682837
$ws.Cells.Item(12, 12).Value = 1.039026450522921
$ws.Cells.Item(12, 13).Value = 1.052152955579728

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.036590806205342
$ws.Cells.Item(13, 4).Value = 1.041379911857992
$ws.Cells.Item(13, 5).Value = 1.035482207955183
$ws.Cells.Item(13, 6).Value = 1.048661326191791
$ws.Cells.Item(13, 9).Value = 1.03432832816017
$ws.Cells.Item(13, 10).Value = 1.043214448275359
$ws.Cells.Item(13, 11).Value = 1.044962310067076
$ws.Cells.Item(13, 12).Value = 1.039086586681659
$ws.Cells.Item(13, 13).Value = 1.052216963848371

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.036857182888419
$ws.Cells.Item(14, 4).Value = 1.041581769064321
$ws.Cells.Item(14, 5).Value = 1.03571167942521
$ws.Cells.Item(14, 6).Value = 1.04890304836171
$ws.Cells.Item(14, 9).Value = 1.034383960509683
$ws.Cells.Item(14, 10).Value = 1.043419239875341
$ws.Cells.Item(14, 11).Value = 1.045132000634343
$ws.Cells.Item(14, 12).Value = 1.03928357416723
$ws.Cells.Item(14, 13).Value = 1.052426636572439

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.037021267447815
$ws.Cells.Item(15, 4).Value = 1.04170610268456
$ws.Cells.Item(15, 5).Value = 1.035853036506352
$ws.Cells.Item(15, 6).Value = 1.049051951690259
$ws.Cells.Item(15, 9).Value = 1.034418188214209
$ws.Cells.Item(15, 10).Value = 1.043545372650118
$ws.Cells.Item(15, 11).Value = 1.045236501618964
$ws.Cells.Item(15, 12).Value = 1.03940490412691
$ws.Cells.Item(15, 13).Value = 1.052555780785692

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.037975723276682
$ws.Cells.Item(16, 4).Value = 1.042429213359699
$ws.Cells.Item(16, 5).Value = 1.036675375414057
$ws.Cells.Item(16, 6).Value = 1.049918185514572
$ws.Cells.Item(16, 9).Value = 1.034616658211118
$ws.Cells.Item(16, 10).Value = 1.044278821156248
$ws.Cells.Item(16, 11).Value = 1.045843968211425
$ws.Cells.Item(16, 12).Value = 1.040110482706685
$ws.Cells.Item(16, 13).Value = 1.053306818311587

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.038573918395485
$ws.Cells.Item(17, 4).Value = 1.04288230883246
$ws.Cells.Item(17, 5).Value = 1.037190842099733
$ws.Cells.Item(17, 6).Value = 1.050461162291122
$ws.Cells.Item(17, 9).Value = 1.034740494256977
$ws.Cells.Item(17, 10).Value = 1.044738283969023
$ws.Cells.Item(17, 11).Value = 1.04622433833332
$ws.Cells.Item(17, 12).Value = 1.040552537556156
$ws.Cells.Item(17, 13).Value = 1.053777367905092

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.038922654671969
$ws.Cells.Item(18, 4).Value = 1.043146416002675
$ws.Cells.Item(18, 5).Value = 1.037491375791688
$ws.Cells.Item(18, 6).Value = 1.050777733569846
$ws.Cells.Item(18, 9).Value = 1.034812488907735
$ws.Cells.Item(18, 10).Value = 1.045006063114693
$ws.Cells.Item(18, 11).Value = 1.046445959296408
$ws.Cells.Item(18, 12).Value = 1.040810189500827
$ws.Cells.Item(18, 13).Value = 1.054051633038856

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.039041534333226
$ws.Cells.Item(19, 4).Value = 1.043236440082559
$ws.Cells.Item(19, 5).Value = 1.037593828283782
$ws.Cells.Item(19, 6).Value = 1.050885653015753
$ws.Cells.Item(19, 9).Value = 1.034836997115402
$ws.Cells.Item(19, 10).Value = 1.045097332148983
$ws.Cells.Item(19, 11).Value = 1.046521485346587
$ws.Cells.Item(19, 12).Value = 1.040898009915098
$ws.Cells.Item(19, 13).Value = 1.054145116886217

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.038509756488728
$ws.Cells.Item(20, 4).Value = 1.042833714196484
$ws.Cells.Item(20, 5).Value = 1.037135550825976
$ws.Cells.Item(20, 6).Value = 1.05040292033456
$ws.Cells.Item(20, 9).Value = 1.034727232331331
$ws.Cells.Item(20, 10).Value = 1.044689010527074
$ws.Cells.Item(20, 11).Value = 1.046183553339897
$ws.Cells.Item(20, 12).Value = 1.040505129072631
$ws.Cells.Item(20, 13).Value = 1.053726902962928

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.036778749716052
$ws.Cells.Item(21, 4).Value = 1.041522334922031
$ws.Cells.Item(21, 5).Value = 1.035644111613808
$ws.Cells.Item(21, 6).Value = 1.048831873402676
$ws.Cells.Item(21, 9).Value = 1.034367588427152
$ws.Cells.Item(21, 10).Value = 1.043358943466766
$ws.Cells.Item(21, 11).Value = 1.045082041620506
$ws.Cells.Item(21, 12).Value = 1.039225574717337
$ws.Cells.Item(21, 13).Value = 1.052364901957233

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.035688705641418
$ws.Cells.Item(22, 4).Value = 1.040696196176042
$ws.Cells.Item(22, 5).Value = 1.034705173374244
$ws.Cells.Item(22, 6).Value = 1.047842803657764
$ws.Cells.Item(22, 9).Value = 1.034139318521245
$ws.Cells.Item(22, 10).Value = 1.042520670499307
$ws.Cells.Item(22, 11).Value = 1.044387256137385
$ws.Cells.Item(22, 12).Value = 1.038419303394552
$ws.Cells.Item(22, 13).Value = 1.051506726034814

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.036266721994798
$ws.Cells.Item(23, 4).Value = 1.041134303662949
$ws.Cells.Item(23, 5).Value = 1.035203039362896
$ws.Cells.Item(23, 6).Value = 1.048367252871477
$ws.Cells.Item(23, 9).Value = 1.034260532958473
$ws.Cells.Item(23, 10).Value = 1.042965247311842
$ws.Cells.Item(23, 11).Value = 1.044755787413525
$ws.Cells.Item(23, 12).Value = 1.038846892419302
$ws.Cells.Item(23, 13).Value = 1.051961837092159

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.03853874905152
$ws.Cells.Item(24, 4).Value = 1.0428556725641
$ws.Cells.Item(24, 5).Value = 1.037160534972424
$ws.Cells.Item(24, 6).Value = 1.050429237794105
$ws.Cells.Item(24, 9).Value = 1.03473322555654
$ws.Cells.Item(24, 10).Value = 1.044711275747578
$ws.Cells.Item(24, 11).Value = 1.046201983072372
$ws.Cells.Item(24, 12).Value = 1.040526551516601
$ws.Cells.Item(24, 13).Value = 1.05374970651162

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.04116824839004
$ws.Cells.Item(25, 4).Value = 1.044846350533815
$ws.Cells.Item(25, 5).Value = 1.039427050552590
$ws.Cells.Item(25, 6).Value = 1.052816671496797
$ws.Cells.Item(25, 9).Value = 1.035272476026071
$ws.Cells.Item(25, 10).Value = 1.046728925872745
$ws.Cells.Item(25, 11).Value = 1.047870718495768
$ws.Cells.Item(25, 12).Value = 1.042468224506278
$ws.Cells.Item(25, 13).Value = 1.055816665287794
